# "len-wei by 5cm also included"
# Populate columns D (length) and E (weight) on the "len-wei Gadget" sheet
# with every 5th row of the existing A/B length-weight table (i.e. the
# same relationship, sampled in 5 cm length increments instead of 1 cm).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("len-wei Gadget")

# Pairs of (length cm, weight kg) taken from column A/B every 5th row,
# starting at A3/B3 (110.5) through A143/B143 (250.5).
$data = @(
  @(110.5, 24.685580000000002),
  @(115.5, 26.608730000000001),
  @(120.5, 28.681699999999999),
  @(125.5, 30.916170000000001),
  @(130.5, 33.324719999999999),
  @(135.5, 35.920909999999999),
  @(140.5, 38.719349999999999),
  @(145.5, 41.735810000000001),
  @(150.5, 44.987270000000002),
  @(155.5, 48.492040000000003),
  @(160.5, 52.269849999999998),
  @(165.5, 56.341970000000003),
  @(170.5, 60.73133),
  @(175.5, 65.462649999999996),
  @(180.5, 70.562569999999994),
  @(185.5, 76.059799999999996),
  @(190.5, 81.985290000000006),
  @(195.5, 88.372420000000005),
  @(200.5, 95.257140000000007),
  @(205.5, 102.67821000000001),
  @(210.5, 110.67744),
  @(215.5, 119.29985000000001),
  @(220.5, 128.59398999999999),
  @(225.5, 138.6122),
  @(230.5, 149.41088999999999),
  @(235.5, 161.05085),
  @(240.5, 173.59764000000001),
  @(245.5, 187.12189000000001),
  @(250.5, 201.69977)
)

$row = 3
foreach ($pair in $data) {
    $ws.Cells.Item($row, 4).Value = $pair[0]
    $ws.Cells.Item($row, 5).Value = $pair[1]
    $row = $row + 1
}

$ws.Range("G13").Select() | Out-Null
